# Adds a new "Player Info" worksheet (with the player's basic bio info) in
# front of the existing "ODI Batting" / "ODI Bowling" sheets, and replaces the
# full scorecard URL stored in the MATCH_CARD_LINK column of those two sheets
# with just the bare numeric match code, renaming that column to MATCH_CODE.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand new "Player Info" sheet before the batting sheet so the
#    final tab order is: Player Info, ODI Batting, ODI Bowling.
#    NOTE: sheet variables captured here are positional - once a sheet is
#    inserted in front of them their "identity" shifts, so we re-look the
#    batting/bowling sheets up *by name* after the insert instead of reusing
#    a stale handle.
# ---------------------------------------------------------------------------
$battingSheetBeforeInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetBeforeInsert)
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Mirror the bold / bordered / centered-top header look used on the other
# two sheets' header rows.
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

$playerInfo.Range("A2").Value = 4930
$playerInfo.Range("B2").Value = "Muralikrishna Prasidh Krishna"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

$playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE and replace each row's full
#    scorecard URL with just the trailing numeric "MatchCode" value, on both
#    the batting and bowling sheets.
# ---------------------------------------------------------------------------
$matchCodes = @(4454, 4456, 4457, 4529, 4533, 4535, 4536, 4609, 4613, 4618, 4621, 4624, 4637, 4640)

$battingSheet.Range("D1").Value = "MATCH_CODE"
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $battingSheet.Range("D$row").Value = $matchCodes[$i]
}

$bowlingSheet.Range("B1").Value = "MATCH_CODE"
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $bowlingSheet.Range("B$row").Value = $matchCodes[$i]
}
